# Bosses.xlsx refactor: drop TradeRowSize / HandSize / Initiative / Cards columns,
# keeping only Id, Name, Health, Manna (Table1 shrinks from A1:H49 to A1:D49).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bosses")

# Swap the contents of columns D (TradeRowSize) and E (Manna) - header + the 6 data
# rows - so "Manna" ends up in column D (its final resting place) and "TradeRowSize"
# ends up in column E, ready to be removed along with the other trailing columns.
$colD = $ws.Range("D1:D7").Value2
$colE = $ws.Range("E1:E7").Value2
$ws.Range("E1:E7").Value2 = $colD
$ws.Range("D1:D7").Value2 = $colE

# Remove the unwanted columns from the table, always deleting the column that is
# currently last (Cards -> Initiative -> HandSize -> TradeRowSize), leaving the
# table as Id, Name, Health, Manna.
$ws.ListObjects.Item("Table1").ListColumns.Item("Cards").Delete()
$ws.ListObjects.Item("Table1").ListColumns.Item("Initiative").Delete()
$ws.ListObjects.Item("Table1").ListColumns.Item("HandSize").Delete()
$ws.ListObjects.Item("Table1").ListColumns.Item("TradeRowSize").Delete()

# Match the saved selection state.
$ws.Range("F14").Select()
